$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1423430998685689
$ws.Range("D2").Value = 0.07158483997343978
$ws.Range("E2").Value = 0.07745764676342759
$ws.Range("F2").Value = 1.882501660108005
$ws.Range("G2").Value = 0.002491788268177655
$ws.Range("K2").Value = 1.272811968774079
$ws.Range("M2").Value = 0.4086679317595383
$ws.Range("N2").Value = 2.257329970245308

$ws.Range("B3").Value = 0.1329860782556551
$ws.Range("D3").Value = 0.07194469834303163
$ws.Range("E3").Value = 0.07184626472631805
$ws.Range("F3").Value = 1.834881906873136
$ws.Range("G3").Value = 0.00249717793146583
$ws.Range("K3").Value = 1.15332223790071
$ws.Range("M3").Value = 0.3728690383117055
$ws.Range("N3").Value = 2.262382466212301

$ws.Range("B4").Value = 0.1273150356605015
$ws.Range("D4").Value = 0.07218726011375765
$ws.Range("E4").Value = 0.06845219826029236
$ws.Range("F4").Value = 1.806765004655333
$ws.Range("G4").Value = 0.002500658947520624
$ws.Range("K4").Value = 1.080647570257838
$ws.Range("M4").Value = 0.3511290626892603
$ws.Range("N4").Value = 2.266127998233017

$ws.Range("B5").Value = 0.1250227743728942
$ws.Range("D5").Value = 0.07229151198261619
$ws.Range("E5").Value = 0.06708177124968628
$ws.Range("F5").Value = 1.795587397134881
$ws.Range("G5").Value = 0.002502120829616528
$ws.Range("K5").Value = 1.051203448059454
$ws.Range("M5").Value = 0.342329374868342
$ws.Range("N5").Value = 2.267815513991565

$ws.Range("B6").Value = 0.1246432798875929
$ws.Range("D6").Value = 0.07230914864330629
$ws.Range("E6").Value = 0.06685497225583958
$ws.Range("F6").Value = 1.793748237733269
$ws.Range("G6").Value = 0.002502366196151109
$ws.Range("K6").Value = 1.04632454650249
$ws.Range("M6").Value = 0.3408717599831945
$ws.Range("N6").Value = 2.268105442861554

$ws.Range("B7").Value = 0.1272840454658848
$ws.Range("D7").Value = 0.07218864423884241
$ws.Range("E7").Value = 0.0684336651493993
$ws.Range("F7").Value = 1.806613127208777
$ws.Range("G7").Value = 0.002500678487278869
$ws.Range("K7").Value = 1.080249786151768
$ws.Range("M7").Value = 0.3510101472946232
$ws.Range("N7").Value = 2.266150104834651

$ws.Range("B8").Value = 0.1391014386017417
$ws.Range("D8").Value = 0.07170441795864235
$ws.Range("E8").Value = 0.07551204577836046
$ws.Range("F8").Value = 1.86584840825688
$ws.Range("G8").Value = 0.002493611072491573
$ws.Range("K8").Value = 1.231466480446613
$ws.Range("M8").Value = 0.3962739049781092
$ws.Range("N8").Value = 2.258938208220485

$ws.Range("B9").Value = 0.1628622019188128
$ws.Range("D9").Value = 0.07092754763001707
$ws.Range("E9").Value = 0.08981122592701496
$ws.Range("F9").Value = 1.991001085977075
$ws.Range("G9").Value = 0.002481107530048487
$ws.Range("K9").Value = 1.533637308493894
$ws.Range("M9").Value = 0.4869956841956906
$ws.Range("N9").Value = 2.249926740487723

$ws.Range("B10").Value = 0.1806765052107977
$ws.Range("D10").Value = 0.07046382200439183
$ws.Range("E10").Value = 0.1005884801463353
$ws.Range("F10").Value = 2.088571880402924
$ws.Range("G10").Value = 0.002472737717014231
$ws.Range("K10").Value = 1.759296755071546
$ws.Range("M10").Value = 0.5549200920212343
$ws.Range("N10").Value = 2.246472598048115

$ws.Range("B11").Value = 0.1888583060809452
$ws.Range("D11").Value = 0.07027649615865883
$ws.Range("E11").Value = 0.1055539900104066
$ws.Range("F11").Value = 2.134211106057762
$ws.Range("G11").Value = 0.002469105259967509
$ws.Range("K11").Value = 1.862798012059102
$ws.Range("M11").Value = 0.5861138132760431
$ws.Range("N11").Value = 2.245597340192361

$ws.Range("B12").Value = 0.1919677145311596
$ws.Range("D12").Value = 0.07020899142452564
$ws.Range("E12").Value = 0.1074436219974331
$ws.Range("F12").Value = 2.15167612170174
$ws.Range("G12").Value = 0.002467754747616251
$ws.Range("K12").Value = 1.902116723378356
$ws.Range("M12").Value = 0.5979696903906131
$ws.Range("N12").Value = 2.245366688098855

$ws.Range("B13").Value = 0.1912975536206147
$ws.Range("D13").Value = 0.0702233766063749
$ws.Range("E13").Value = 0.1070362381088259
$ws.Range("F13").Value = 2.147906569968455
$ws.Range("G13").Value = 0.002468044494240845
$ws.Range("K13").Value = 1.89364312394099
$ws.Range("M13").Value = 0.5954143625636874
$ws.Range("N13").Value = 2.245411869306182

$ws.Range("B14").Value = 0.1891138961956642
$ws.Range("D14").Value = 0.07027087353190353
$ws.Range("E14").Value = 0.1057092632668102
$ws.Range("F14").Value = 2.13564429483381
$ws.Range("G14").Value = 0.002468993651921829
$ws.Range("K14").Value = 1.866030262082461
$ws.Range("M14").Value = 0.5870883252194119
$ws.Range("N14").Value = 2.245576340311672

$ws.Range("B15").Value = 0.1877777912290952
$ws.Range("D15").Value = 0.07030041468135906
$ws.Range("E15").Value = 0.1048976718889207
$ws.Range("F15").Value = 2.128157119276352
$ws.Range("G15").Value = 0.002469578292782454
$ws.Range("K15").Value = 1.849132963917896
$ws.Range("M15").Value = 0.5819940848257374
$ws.Range("N15").Value = 2.245690230140838

$ws.Range("B16").Value = 0.1801433680091691
$ws.Range("D16").Value = 0.07047654249754487
$ws.Range("E16").Value = 0.1002652597917404
$ws.Range("F16").Value = 2.085614661034668
$ws.Range("G16").Value = 0.002472978616029927
$ws.Range("K16").Value = 1.752550028495932
$ws.Range("M16").Value = 0.5528875328860607
$ws.Range("N16").Value = 2.246543863943302

$ws.Range("B17").Value = 0.1754798159803954
$ws.Range("D17").Value = 0.07059066631250488
$ws.Range("E17").Value = 0.09743971371644733
$ws.Range("F17").Value = 2.059838875544301
$ws.Range("G17").Value = 0.002475109324677173
$ws.Range("K17").Value = 1.693518844784421
$ws.Range("M17").Value = 0.5351078523468829
$ws.Range("N17").Value = 2.24724630908625

$ws.Range("B18").Value = 0.1728048089989898
$ws.Range("D18").Value = 0.07065852894271885
$ws.Range("E18").Value = 0.09582043748735458
$ws.Range("F18").Value = 2.045131207386731
$ws.Range("G18").Value = 0.002476351333812429
$ws.Range("K18").Value = 1.659645262049935
$ws.Range("M18").Value = 0.5249091075235839
$ws.Range("N18").Value = 2.247715810068655

$ws.Range("B19").Value = 0.1719003614103656
$ws.Range("D19").Value = 0.07068188660026919
$ws.Range("E19").Value = 0.09527318442842159
$ws.Range("F19").Value = 2.040171625363627
$ws.Range("G19").Value = 0.002476774692063977
$ws.Range("K19").Value = 1.648189824325925
$ws.Range("M19").Value = 0.5214607036841699
$ws.Range("N19").Value = 2.247885998974027

$ws.Range("B20").Value = 0.1759754994505442
$ws.Range("D20").Value = 0.07057828746077632
$ws.Range("E20").Value = 0.09773988510513476
$ws.Range("F20").Value = 2.062570536200468
$ws.Range("G20").Value = 0.002474880802195848
$ws.Range("K20").Value = 1.699794554936489
$ws.Range("M20").Value = 0.5369976591303072
$ws.Range("N20").Value = 2.247164751359534

$ws.Range("B21").Value = 0.1897549876326536
$ws.Range("D21").Value = 0.07025682912848552
$ws.Range("E21").Value = 0.1060987734331675
$ws.Range("F21").Value = 2.139241053925417
$ws.Range("G21").Value = 0.002468714183403628
$ws.Range("K21").Value = 1.874137410311391
$ws.Range("M21").Value = 0.5895326947300674
$ws.Range("N21").Value = 2.245525290191239

$ws.Range("B22").Value = 0.1988255660941007
$ws.Range("D22").Value = 0.07006675533205708
$ws.Range("E22").Value = 0.1116161272365233
$ws.Range("F22").Value = 2.190414036859181
$ws.Range("G22").Value = 0.002464829710576583
$ws.Range("K22").Value = 1.988810823853953
$ws.Range("M22").Value = 0.6241214209634052
$ws.Range("N22").Value = 2.245041614181943

$ws.Range("B23").Value = 0.1939785154334999
$ws.Range("D23").Value = 0.07016635816265904
$ws.Range("E23").Value = 0.1086663533764636
$ws.Range("F23").Value = 2.163003947204231
$ws.Range("G23").Value = 0.002466889637043022
$ws.Range("K23").Value = 1.92753959118528
$ws.Range("M23").Value = 0.6056371423594413
$ws.Range("N23").Value = 2.245245743627862

$ws.Range("B24").Value = 0.1757513818508443
$ws.Range("D24").Value = 0.07058387692972445
$ws.Range("E24").Value = 0.09760416154311713
$ws.Range("F24").Value = 2.061335206296008
$ws.Range("G24").Value = 0.002474984064135821
$ws.Range("K24").Value = 1.696957106547529
$ws.Range("M24").Value = 0.5361432058344775
$ws.Range("N24").Value = 2.247201419156937

$ws.Range("B25").Value = 0.1563715120775413
$ws.Range("D25").Value = 0.07111904388905899
$ws.Range("E25").Value = 0.08589642052932334
$ws.Range("F25").Value = 1.956167005862312
$ws.Range("G25").Value = 0.002484345963948401
$ws.Range("K25").Value = 1.451266338866901
$ws.Range("M25").Value = 0.4622354030037314
$ws.Range("N25").Value = 2.251811314629492

